# Revert capacity chart to show kilowatts (not watts) on the y-axis.
#  - The underlying "Solar" figures were entered in watts (5000, 3770,
#    10150); they should be kilowatts (5, 3.77, 10.15).
#  - The shared number format used by the data table needs a decimal so the
#    fractional kW values (3.77, 10.15) aren't rounded away.
#  - The value-axis title/number format on the chart need to go back to a
#    plain "Kilowatts (kW)" label with a plain thousands-separated number
#    format instead of the "###K" watts-oriented one.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Worksheet data: convert the three Solar (column E) entries from watts
#     to kilowatts ---------------------------------------------------------
$ws.Range("E22").Value = 5
$ws.Range("E25").Value = 3.77
$ws.Range("E26").Value = 10.15

# --- Number format for the data table: allow one decimal place so values
#     like 3.77 / 10.15 kW display correctly instead of rounding to whole
#     numbers -----------------------------------------------------------
$ws.Range("B2:G26").NumberFormat = "#,##0.0"

# --- Chart: restore the kilowatt axis title + plain number format --------
$co = $ws.ChartObjects(1)
$chart = $co.Chart
$valueAxis = $chart.Axes(2)  # xlValue
$valueAxis.AxisTitle.Text = "Kilowatts (kW)"
$valueAxis.TickLabels.NumberFormat = "#,##0"
